$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $find = $d.Content.Find
    $found = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output ("NOT FOUND: " + $old)
        return
    }
    $rng = $find.Parent
    $rng.Text = $new
}

function InsertAfter-Text($anchorOld, $newText) {
    # Finds anchorOld, collapses to its end, and inserts newText right after it.
    $find = $d.Content.Find
    $found = $find.Execute($anchorOld, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output ("NOT FOUND (anchor): " + $anchorOld)
        return
    }
    $rng = $find.Parent
    $rng.Collapse(0)
    $rng.InsertAfter($newText)
}

# ---------------------------------------------------------------------------
# 1. Title
# ---------------------------------------------------------------------------
Replace-Text "Quantum Entanglement: Unveiling the Mysteries of Interconnectedness" "Mathematics: The Universal Language of Science and Logic"

# ---------------------------------------------------------------------------
# 2. Author name
# ---------------------------------------------------------------------------
Replace-Text " Samuel Davies" " Jasper Lancaster"

# ---------------------------------------------------------------------------
# 3. Email line -> collapses to a single run "yourvalidname"
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$s3 = $p3.Range.Start
$e3 = $p3.Range.End
$r3 = $d.Range($s3, $e3)
$r3.Text = "yourvalidname"

# ---------------------------------------------------------------------------
# 4. Body paragraph (paragraph 5)
# ---------------------------------------------------------------------------
Replace-Text "Amidst the perplexing wonders of the quantum realm exists a profound phenomenon challenging our understanding of reality - quantum entanglement" "Mathematics, a subject that has fascinated and challenged minds for centuries, unveils the mysteries of the universe through the intricate tapestry of numbers, equations, and geometric patterns"

Replace-Text " This enigmatic connection between particles, regardless of their distance, has captured the imagination of scientists, philosophers, and artists alike" " Mathematics provides a universal language that transcends cultures, enabling us to understand the cosmos, unravel the enigmas of nature, and harness the power of logic to solve complex problems"

Replace-Text " In this essay, we embark on a journey to unravel the mysteries of quantum entanglement, exploring its implications for our comprehension of the universe and delving into the potential applications that may revolutionize various fields" " The field of mathematics is a symphony of abstract concepts, where symbols dance in harmony, revealing the underlying order and beauty of our world"

Replace-Text "In 1935, Albert Einstein, Boris Podolsky, and Nathan Rosen introduced the concept of quantum entanglement through their famous thought experiment known as the EPR paradox" "Through the exploration of mathematical concepts, we unlock the secrets of nature's blueprint"

Replace-Text " Their proposal demonstrated that two particles, once entangled, remain interconnected regardless of the distance separating them" " The Fibonacci sequence, found in the spirals of seashells or the patterns of plant growth, exemplifies the intricate relationship between numbers and biological structures"

Replace-Text " This relationship transcends the constraints of space and time, allowing one particle to instantaneously influence the other, even across vast cosmological distances" " The elegance of geometric shapes, such as fractals, reflects the self-similarity found in everything from snowflakes to coastlines"

# New sentence + period inserted right after "...coastlines." and before the next <w:br/>
InsertAfter-Text "The elegance of geometric shapes, such as fractals, reflects the self-similarity found in everything from snowflakes to coastlines." " These patterns underscore the profound interconnectedness of all things and provide a glimpse into the underlying mathematical principles that govern our universe."

Replace-Text "Moreover, quantum entanglement defies classical intuition" "Mathematics isn't merely a collection of abstract theories; it's a powerful tool with practical applications in every field imaginable"

Replace-Text " When entangled particles are measured, their properties, such as spin or polarization, are correlated in a way that cannot be explained by classical physics" " It empowers engineers to design structures that withstand earthquakes, enables us to predict weather patterns, and makes it possible to develop new medical treatments and technologies"

Replace-Text " This non-locality, as it is known, challenges our conventional notions of causality and raises fundamental questions about the nature of reality itself" " From the economy to finance, from computer science to data analysis, and even in music and art, the profound influence of mathematics is undeniable"

# New sentence + period inserted at the very end of the paragraph
InsertAfter-Text "From the economy to finance, from computer science to data analysis, and even in music and art, the profound influence of mathematics is undeniable." " It's a subject that touches every aspect of our lives, shaping our understanding of the world and guiding us towards a future filled with infinite possibilities."

# ---------------------------------------------------------------------------
# 5. Summary body (paragraph 7)
# ---------------------------------------------------------------------------
Replace-Text "Quantum entanglement, an awe-inspiring phenomenon, offers a glimpse into the uncharted territory of the quantum world" "Mathematics serves as an essential tool for understanding the intricacies of our universe, providing a lens through which we unlock the mysteries of science and logic"

Replace-Text " Its non-local nature challenges our fundamental understanding of reality, while its potential applications hold promise for transformative technologies" " Its abstract concepts find practical applications in diverse fields, empowering engineers, scientists, musicians, artists, and countless other professionals to innovate and drive progress"

Replace-Text " From quantum computing to secure communication, entanglement-based technologies may revolutionize numerous fields" " Mathematics unveils the interconnectedness of all things, from the Fibonacci sequence found in nature to the intricate "

Replace-Text " Though much remains unknown, continued exploration of quantum entanglement promises to deepen our comprehension of the universe and expand the boundaries of human knowledge" " It's a subject that permeates our existence, shaping our understanding of the cosmos and enabling us to chart a course toward a future of endless possibilities"

# Insert the "patterns in art and music" run (with a lastRenderedPageBreak) right before the
# "." that follows "...to the intricate "
$find = $d.Content.Find
$found = $find.Execute("to the intricate ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng = $find.Parent
    $rng.Collapse(0)
    $rng.InsertAfter("patterns in art and music")
}

# ---------------------------------------------------------------------------
# 6. New empty paragraph appended at the very end of the document
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 7. Global font rename: TimesNewToman -> Times New Roman
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.Trim().Length -gt 0) {
        $p.Range.Font.Name = "Times New Roman"
    }
}
